# Update automatic: dades i banners [2026-02-13 21:20]
# Refreshes the per-station extraction timestamps and the latest observation
# values (snow depth, humidity, precipitation, pressure, wind gust, temps...)
# in the "Dades_Meteo" sheet to match the 21:18-21:20 meteo.cat pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('E2').Value = '2026-02-13 21:18:34'
$ws.Range('G2').Value = '177 cm'
$ws.Range('I2').Value = '3.0 mm'
$ws.Range('O2').Value = '-0.5 °C'

# Row 3
$ws.Range('E3').Value = '2026-02-13 21:18:37'
$ws.Range('I3').Value = '6.9 mm'

# Row 4
$ws.Range('E4').Value = '2026-02-13 21:18:39'
$ws.Range('J4').Value = '993.8 hPa'

# Row 5
$ws.Range('E5').Value = '2026-02-13 21:18:42'
$ws.Range('I5').Value = '2.7 mm'

# Row 6
$ws.Range('E6').Value = '2026-02-13 21:18:44'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '78%'
$ws.Range('J6').Value = '993.8 hPa'

# Row 7
$ws.Range('E7').Value = '2026-02-13 21:18:46'
$ws.Range('J7').Value = '994.2 hPa'
$ws.Range('O7').Value = '12.7 °C'

# Row 8
$ws.Range('E8').Value = '2026-02-13 21:18:48'
$ws.Range('J8').Value = '994.1 hPa'
$ws.Range('L8').Value = '49.7 km/h - 309º 20:39 TU'
$ws.Range('N8').Value = '6.5 °C 20:35 TU'
$ws.Range('O8').Value = '9.1 °C'

# Row 9
$ws.Range('E9').Value = '2026-02-13 21:18:51'

# Row 10
$ws.Range('E10').Value = '2026-02-13 21:18:53'

# Row 11
$ws.Range('E11').Value = '2026-02-13 21:18:56'
$ws.Range('H11').NumberFormat = '@'
$ws.Range('H11').Value = '93%'

# Row 12
$ws.Range('E12').Value = '2026-02-13 21:18:58'

# Row 13
$ws.Range('E13').Value = '2026-02-13 21:19:00'
$ws.Range('J13').Value = '996.9 hPa'
$ws.Range('O13').Value = '0.6 °C'

# Row 14
$ws.Range('E14').Value = '2026-02-13 21:19:03'
$ws.Range('L14').Value = '33.8 km/h - 319º 20:52 TU'

# Row 15
$ws.Range('E15').Value = '2026-02-13 21:19:05'
$ws.Range('H15').NumberFormat = '@'
$ws.Range('H15').Value = '76%'

# Row 16
$ws.Range('E16').Value = '2026-02-13 21:19:07'
$ws.Range('I16').Value = '13.6 mm'

# Row 17
$ws.Range('E17').Value = '2026-02-13 21:19:10'
$ws.Range('G17').Value = '2 cm'

# Row 18
$ws.Range('E18').Value = '2026-02-13 21:19:12'
$ws.Range('H18').NumberFormat = '@'
$ws.Range('H18').Value = '84%'
$ws.Range('J18').Value = '994.0 hPa'

# Row 19
$ws.Range('E19').Value = '2026-02-13 21:19:15'

# Row 20
$ws.Range('E20').Value = '2026-02-13 21:19:17'
$ws.Range('I20').Value = '24.1 mm'

# Row 21
$ws.Range('E21').Value = '2026-02-13 21:19:20'
$ws.Range('J21').Value = '997.0 hPa'
$ws.Range('N21').Value = '-0.2 °C 20:54 TU'

# Row 22
$ws.Range('E22').Value = '2026-02-13 21:19:22'

# Row 23
$ws.Range('E23').Value = '2026-02-13 21:19:24'
$ws.Range('I23').Value = '11.5 mm'

# Row 24
$ws.Range('E24').Value = '2026-02-13 21:19:27'
$ws.Range('H24').NumberFormat = '@'
$ws.Range('H24').Value = '95%'
$ws.Range('J24').Value = '995.1 hPa'
$ws.Range('L24').Value = '60.8 km/h - 291º 20:36 TU'
$ws.Range('O24').Value = '7.2 °C'

# Row 25
$ws.Range('E25').Value = '2026-02-13 21:19:29'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '80%'
$ws.Range('I25').Value = '9.4 mm'

# Row 26
$ws.Range('E26').Value = '2026-02-13 21:19:32'

# Row 27
$ws.Range('E27').Value = '2026-02-13 21:19:34'

# Row 28
$ws.Range('E28').Value = '2026-02-13 21:19:37'
$ws.Range('H28').NumberFormat = '@'
$ws.Range('H28').Value = '81%'
$ws.Range('I28').Value = '6.6 mm'
$ws.Range('J28').Value = '994.3 hPa'

# Row 29
$ws.Range('E29').Value = '2026-02-13 21:19:39'

# Row 30
$ws.Range('E30').Value = '2026-02-13 21:19:42'
$ws.Range('J30').Value = '993.7 hPa'
$ws.Range('O30').Value = '9.3 °C'

# Row 31
$ws.Range('E31').Value = '2026-02-13 21:19:44'
$ws.Range('J31').Value = '992.7 hPa'
$ws.Range('O31').Value = '10.2 °C'

# Row 32
$ws.Range('E32').Value = '2026-02-13 21:19:47'
$ws.Range('H32').NumberFormat = '@'
$ws.Range('H32').Value = '91%'

# Row 33
$ws.Range('E33').Value = '2026-02-13 21:19:49'
$ws.Range('J33').Value = '995.8 hPa'

# Row 34
$ws.Range('E34').Value = '2026-02-13 21:19:52'

# Row 35
$ws.Range('E35').Value = '2026-02-13 21:19:54'
$ws.Range('H35').NumberFormat = '@'
$ws.Range('H35').Value = '77%'
$ws.Range('I35').Value = '8.7 mm'
$ws.Range('J35').Value = '995.1 hPa'
$ws.Range('O35').Value = '5.9 °C'

# Row 36
$ws.Range('E36').Value = '2026-02-13 21:19:56'
$ws.Range('J36').Value = '993.9 hPa'

# Row 37
$ws.Range('E37').Value = '2026-02-13 21:19:59'
$ws.Range('J37').Value = '995.8 hPa'

# Row 38
$ws.Range('E38').Value = '2026-02-13 21:20:01'
$ws.Range('H38').NumberFormat = '@'
$ws.Range('H38').Value = '78%'

# Row 39
$ws.Range('E39').Value = '2026-02-13 21:20:04'
$ws.Range('I39').Value = '19.6 mm'

# Row 40
$ws.Range('E40').Value = '2026-02-13 21:20:06'
$ws.Range('J40').Value = '997.5 hPa'

# Row 41
$ws.Range('E41').Value = '2026-02-13 21:20:08'
$ws.Range('J41').Value = '994.4 hPa'
$ws.Range('L41').Value = '61.9 km/h - 283º 20:47 TU'
$ws.Range('O41').Value = '12.1 °C'

# Row 42
$ws.Range('E42').Value = '2026-02-13 21:20:11'

# Row 43
$ws.Range('E43').Value = '2026-02-13 21:20:13'

# Row 44
$ws.Range('E44').Value = '2026-02-13 21:20:16'
$ws.Range('I44').Value = '9.4 mm'

# Row 45
$ws.Range('E45').Value = '2026-02-13 21:20:18'
$ws.Range('H45').NumberFormat = '@'
$ws.Range('H45').Value = '66%'
$ws.Range('I45').Value = '1.7 mm'
$ws.Range('J45').Value = '993.1 hPa'
$ws.Range('O45').Value = '5.5 °C'

# Row 46
$ws.Range('E46').Value = '2026-02-13 21:20:21'
$ws.Range('J46').Value = '995.2 hPa'
$ws.Range('L46').Value = '51.8 km/h - 332º 20:53 TU'
$ws.Range('M46').Value = '12.5 °C 20:59 TU'
$ws.Range('O46').Value = '9.1 °C'
